# The deck's live design theme ("Integral", stored in ppt/theme/theme2.xml and
# referenced by the slide master / presentation) is being swapped for the
# plain default "Office Theme" palette (currently stored in ppt/theme/theme1.xml,
# which is only wired to the notes master).
#
# The font scheme and format scheme are already identical between the two
# theme parts, so only the 12-colour theme colour scheme actually changes.
# We drive this through the PowerPoint theme-colour-scheme object, which is
# exposed per slide and resolves to the presentation's single live theme part.

function ToCOMRGB($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# Target palette: the default Office theme colours (index -> hex), in the
# standard ThemeColorScheme order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $officeColors.Count; $i++) {
    $tcs.Item($i).RGB = ToCOMRGB($officeColors[$i - 1])
}
